$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4697517454624176
$ws.Range("B1").Value = 0.592220664024353
$ws.Range("C1").Value = 0.8604615926742554
$ws.Range("D1").Value = 3.788425207138062
$ws.Range("E1").Value = 5.655208587646484
